# ---------------------------------------------------------------------------
# "added another negative case"
#
# - Renames Sheet1..Sheet4 to TestCase1..TestCase4.
# - Inserts a new "Enter valid email" verification step into TestCase1
#   (the successful-registration case) and TestCase3 (the incorrect-data
#   case), pushing the remaining steps down and renumbering them.
# - Duplicates TestCase4 into a new TestCase5 sheet and edits it into a
#   brand new negative test case covering a taken e-mail address.
# - Leaves TestCase2 essentially untouched and makes TestCase5 the active
#   (selected) sheet, like TestCase4 used to be.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. rename the four existing sheets -----------------------------------
$wb.Worksheets.Item(1).Name = "TestCase1"
$wb.Worksheets.Item(2).Name = "TestCase2"
$wb.Worksheets.Item(3).Name = "TestCase3"
$wb.Worksheets.Item(4).Name = "TestCase4"

$ws1 = $wb.Worksheets.Item("TestCase1")
$ws2 = $wb.Worksheets.Item("TestCase2")
$ws3 = $wb.Worksheets.Item("TestCase3")
$ws4 = $wb.Worksheets.Item("TestCase4")

# --- 2. TestCase1: insert a new step verifying the e-mail field -----------
$ws1.Rows.Item(13).Insert()
$ws1.Cells.Item(13, 1).Value = 3
$ws1.Cells.Item(13, 2).Value = "Enter valid email"
$ws1.Cells.Item(13, 3).Value = "pesho@pesho.com"
$ws1.Cells.Item(13, 4).Value = "A green tick should show, denoting that the mail is unique"
$ws1.Hyperlinks.Add($ws1.Range("C13"), "mailto:pesho@pesho.com") | Out-Null

# renumber the steps that got pushed down
$ws1.Cells.Item(14, 1).Value = 4
$ws1.Cells.Item(15, 1).Value = 5
$ws1.Cells.Item(16, 1).Value = 6

$ws1.Range("D17").Select()

# --- 3. TestCase3: same kind of new step -----------------------------------
$ws3.Rows.Item(13).Insert()
$ws3.Cells.Item(13, 1).Value = 3
$ws3.Cells.Item(13, 2).Value = "Enter valid mail"
$ws3.Cells.Item(13, 3).Value = "gosho@gosho.com"
$ws3.Cells.Item(13, 4).Value = "A green tick should show, denoting that the email is unique"
$ws3.Hyperlinks.Add($ws3.Range("C13"), "mailto:gosho@gosho.com") | Out-Null

$ws3.Cells.Item(14, 1).Value = 4
$ws3.Cells.Item(15, 1).Value = 5
$ws3.Cells.Item(16, 1).Value = 6

$ws3.Range("C13").Select()

# --- 4. TestCase4: no content change, just leaves the "active sheet" spot -
$ws4.Range("B18").Select()

# --- 5. TestCase5: duplicate of TestCase4, turned into a new negative case
$ws4.Copy([Type]::Missing, $ws4)
$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "TestCase5"

$ws5.Cells.Item(12, 2).Value = "Enter a username"
$ws5.Cells.Item(12, 3).Value = "GoGO#67619"
$ws5.Cells.Item(12, 4).ClearContents()

$ws5.Cells.Item(13, 1).Value = 3
$ws5.Cells.Item(13, 2).Value = "Enter email"
$ws5.Cells.Item(13, 3).Value = "gosho@gosho.com"
$ws5.Cells.Item(13, 4).Value = "A red sign should be shown, along with a writing saying that the email is taken"
$ws5.Hyperlinks.Add($ws5.Range("C13"), "mailto:gosho@gosho.com") | Out-Null
$ws5.Range("D13").WrapText = $true
$ws5.Rows.Item(13).RowHeight = 30

$ws5.Range("D16").Select()
$ws5.Activate()

Write-Output "done"
